$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1535.6830860763801
$ws.Range("H2").Value = 15225.7511101741
$ws.Range("I2").Value = 23522.883315536299
$ws.Range("J2").Value = 19357.243296258199
$ws.Range("K2").Value = 18464.9123915987
$ws.Range("R2").Value = 0.69899999999999995
$ws.Range("S2").Value = 0.88900000000000001
$ws.Range("T2").Value = 0.84899999999999998
$ws.Range("U2").Value = 0.77700000000000002
$ws.Range("AA2").Value = 23625.893631944298
$ws.Range("AB2").Value = 21782.190429433598
$ws.Range("AC2").Value = 26459.936237948601
$ws.Range("AD2").Value = 22800.050996770598
$ws.Range("AE2").Value = 23764.3660123523
$ws.Range("AM2").Value = 0.00000000039402810099958599
$ws.Range("AN2").Value = 77.392788372883302
$ws.Range("AO2").Value = 130.70909123094799
$ws.Range("AR2").Value = 78106.473199644097
$ws.Range("AS2").Value = 50697.273153351402
$ws.Range("AT2").Value = 184.54846031123901
$ws.Range("G3").Value = 14644.6858046045
$ws.Range("H3").Value = 27407.792237613001
$ws.Range("I3").Value = 23766.275375901299
$ws.Range("J3").Value = 19277.522172099401
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 0.91500000000000004
$ws.Range("T3").Value = 0.80900000000000005
$ws.Range("AA3").Value = 23135.364620228302
$ws.Range("AB3").Value = 27407.792237613001
$ws.Range("AC3").Value = 25974.0714490725
$ws.Range("AD3").Value = 23828.828395672899
$ws.Range("AL3").Value = 29.939842693673999
$ws.Range("AM3").Value = 1.4722042540484901
$ws.Range("AN3").Value = 32.395950729590403
$ws.Range("AR3").Value = 85096.275590217905
$ws.Range("AS3").Value = 57777.725226495997
$ws.Range("AT3").Value = 275.19814288144102
$ws.Range("F4").Value = 13557.8254831629
$ws.Range("G4").Value = 27457.739379734499
$ws.Range("H4").Value = 28026.869818484
$ws.Range("I4").Value = 22594.787147544299
$ws.Range("S4").Value = 0.81200000000000006
$ws.Range("Z4").Value = 27725.6144849958
$ws.Range("AA4").Value = 27457.739379734499
$ws.Range("AB4").Value = 28139.4275285984
$ws.Range("AC4").Value = 27826.092546236901
$ws.Range("AL4").Value = -0.00000000022786858789551101
$ws.Range("AM4").Value = 69.782604835440395
$ws.Range("AR4").Value = 91637.221828923794
$ws.Range("AS4").Value = 64381.757472895697
$ws.Range("AT4").Value = 338.284150573413
$ws.Range("E5").Value = 13547.8118764602
$ws.Range("F5").Value = 25430.750109982499
$ws.Range("G5").Value = 32678.171766084299
$ws.Range("H5").Value = 24076.056006483501
$ws.Range("Y5").Value = 27705.136761677299
$ws.Range("Z5").Value = 25430.750109982499
$ws.Range("AA5").Value = 32678.171766084299
$ws.Range("AB5").Value = 24567.404088248499
$ws.Range("AJ5").Value = 0.17778624373571
$ws.Range("AK5").Value = 0.0000000000199341603431291
$ws.Range("AL5").Value = 7.9164850868996899
$ws.Range("AM5").Value = -0.00000000000000016177825619341999
$ws.Range("AR5").Value = 95732.789759008796
$ws.Range("AS5").Value = 68449.081454617306
$ws.Range("AT5").Value = 310.04020221046397
$ws.Range("H6").Value = 238
$ws.Range("I6").Value = 1568.10975872376
$ws.Range("J6").Value = 9869.8458710506002
$ws.Range("K6").Value = 25990.539275126699
$ws.Range("L6").Value = 19360.282094247599
$ws.Range("M6").Value = 15925.687364056201
$ws.Range("T6").Value = 0.0070000000000000001
$ws.Range("U6").Value = 0.17399999999999999
$ws.Range("V6").Value = 0.96199999999999997
$ws.Range("W6").Value = 0.98399999999999999
$ws.Range("X6").Value = 0.79100000000000004
$ws.Range("Y6").Value = 0.63300000000000001
$ws.Range("AG6").Value = 9012.1250501365994
$ws.Range("AH6").Value = 10259.7150426721
$ws.Range("AI6").Value = 26413.149669844199
$ws.Range("AJ6").Value = 24475.704289009602
$ws.Range("AK6").Value = 25159.063753319599
$ws.Range("AU6").Value = 23.8237949470785
$ws.Range("AV6").Value = 6.7409482626257597
$ws.Range("AW6").Value = 107.69572988842501
$ws.Range("AZ6").Value = 72952.464363205101
$ws.Range("BA6").Value = 49349.571813202703
$ws.Range("BB6").Value = 255.59195194692299
$ws.Range("H7").Value = 348.574876358783
$ws.Range("I7").Value = 5958.6700849487697
$ws.Range("J7").Value = 30565.014191755901
$ws.Range("K7").Value = 25217.5469133732
$ws.Range("L7").Value = 18557.016376922598
$ws.Range("W7").Value = 0.94699999999999995
$ws.Range("X7").Value = 0.79900000000000004
$ws.Range("AF7").Value = 17428.743817939099
$ws.Range("AG7").Value = 5976.5998846025796
$ws.Range("AH7").Value = 30626.2667252063
$ws.Range("AI7").Value = 26628.877416444899
$ws.Range("AJ7").Value = 23225.302098776701
$ws.Range("AU7").Value = -0.0000000000031338774934949201
$ws.Range("AV7").Value = 10.788178300779
$ws.Range("AW7").Value = -0.0000000000000248717004138803
$ws.Range("AZ7").Value = 80646.822443359095
$ws.Range("BA7").Value = 57222.045606976899
$ws.Range("BB7").Value = 433.707665566819
$ws.Range("H8").Value = 4971.2507151616101
$ws.Range("I8").Value = 33317.8047554105
$ws.Range("J8").Value = 28677.1331406764
$ws.Range("K8").Value = 22152.563662696401
$ws.Range("V8").Value = 0.99299999999999999
$ws.Range("W8").Value = 0.79200000000000004
$ws.Range("AF8").Value = 4971.2507151616101
$ws.Range("AG8").Value = 33317.8047554105
$ws.Range("AH8").Value = 28879.288157780898
$ws.Range("AI8").Value = 27970.408665020699
$ws.Range("AS8").Value = 0
$ws.Range("AT8").Value = 0
$ws.Range("AU8").Value = 2.25115153486886
$ws.Range("AV8").Value = 0
$ws.Range("AW8").Value = 0
$ws.Range("AZ8").Value = 89118.752273944003
$ws.Range("BA8").Value = 65570.485774551998
$ws.Range("BB8").Value = 310.21800255608701
$ws.Range("G9").Value = 4378.82726051386
$ws.Range("H9").Value = 34000
$ws.Range("I9").Value = 30602.999020662999
$ws.Range("J9").Value = 24562.8772041907
$ws.Range("V9").Value = 0.92700000000000005
$ws.Range("AE9").Value = 4378.82726051386
$ws.Range("AF9").Value = 34000
$ws.Range("AG9").Value = 30602.999020662999
$ws.Range("AH9").Value = 26497.170662557401
$ws.Range("AR9").Value = 0
$ws.Range("AS9").Value = 0
$ws.Range("AT9").Value = -0.000000000035644516436120698
$ws.Range("AZ9").Value = 93544.703485366903
$ws.Range("BA9").Value = 69959.401329425702
$ws.Range("BB9").Value = 273.18234600724799

# Apply number-format changes (scientific notation for tiny residuals)
$ws.Range("AM2").NumberFormat = "0.00E+00"
$ws.Range("AL4").NumberFormat = "0.00E+00"
$ws.Range("AK5").NumberFormat = "0.00E+00"
$ws.Range("AM5").NumberFormat = "0.00E+00"
$ws.Range("AU7").NumberFormat = "0.00E+00"
$ws.Range("AW7").NumberFormat = "0.00E+00"
$ws.Range("AT9").NumberFormat = "0.00E+00"

# Reset cells back to default/general style
$ws.Range("AU6").Style = "Normal"
$ws.Range("AV6").Style = "Normal"
$ws.Range("AS8").Style = "Normal"
$ws.Range("AT8").Style = "Normal"
$ws.Range("AV8").Style = "Normal"
$ws.Range("AW8").Style = "Normal"
$ws.Range("AR9").Style = "Normal"

# Update selection to match author's final cursor position
$ws.Range("F12").Select()
